$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# Add three new rows (33-35) describing the PERSON object type's default
# Access Control List participants, mirroring the existing ORGANIZATION rows
# (30-32) immediately above them.
# ---------------------------------------------------------------------------

# Copy the cell formatting (styles) from the ORGANIZATION rows down onto the
# new rows first, so the new cells inherit the same look (borders / number
# format / etc.) as the rest of the ACL table.
$ws.Range("B30:H30").Copy() | Out-Null
$ws.Range("B33:H33").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B31:H31").Copy() | Out-Null
$ws.Range("B34:H34").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B31:H31").Copy() | Out-Null
$ws.Range("B35:H35").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("D32").Copy() | Out-Null
$ws.Range("D35").PasteSpecial(-4122) | Out-Null       # xlPasteFormats

$excel.CutCopyMode = $false

# Populate the Object Type (column C) cells first so "PERSON" becomes the
# next new shared string, then fill in the rest of each row left to right.

# Row 33 - Person - Default access
$ws.Range("C33").Value = "PERSON"
$ws.Range("B33").Value = "Person – Default access"
$ws.Range("D33").Value = "participants.?[participantType == '*'].isEmpty()"
$ws.Range("G33").Value = "*, *"

# Row 34 - Person - Default owner
$ws.Range("C34").Value = "PERSON"
$ws.Range("B34").Value = "Person – Default owner"
$ws.Range("D34").Value = "participants.?[participantType == 'owner'].isEmpty()"
$ws.Range("H34").Value = "owner, creator"

# Row 35 - Person - Default group
$ws.Range("C35").Value = "PERSON"
$ws.Range("B35").Value = "Person – Default group"
$ws.Range("D35").Value = "participants.?[participantType == 'owning group'].isEmpty()"
$ws.Range("G35").Value = "owning group, ACM_ADMINISTRATOR_DEV"

# ---------------------------------------------------------------------------
# Update the sheet view to match where the author left the cursor after
# adding the new rows.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 27
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B36").Select()
